# Updates Sheet1 rows 2-10 with the new TPM-derived NATMI ligand-receptor
# statistics for the Vcan-Tlr2 pair across all Sending x Target cluster
# combinations (ECs, FAPs, MuSCs) x (ECs, FAPs, MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("ECs", "Vcan", "Tlr2", "ECs", 3, 1, 0.2364713333333333, 0.709414, 0.002249544876489787, 0.002249544876489787, 2, 0.6666666666666666, 0.493155, 1.479465, 0.2262533155038342, 0.2262533155038342, 0.11661702039, 1.04955318351, 0.0005089669866804777, 0.0005089669866804777),
    @("ECs", "Vcan", "Tlr2", "FAPs", 3, 1, 0.2364713333333333, 0.709414, 0.002249544876489787, 0.002249544876489787, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.7564224123165462, 0.7564224123165462, 0.3898803767102222, 3.508923390392, 0.001701606162088732, 0.001701606162088732),
    @("ECs", "Vcan", "Tlr2", "MuSCs", 3, 1, 0.2364713333333333, 0.709414, 0.002249544876489787, 0.002249544876489787, 1, 0.3333333333333333, 0.037761, 0.113283, 0.01732427217961956, 0.01732427217961956, 0.008929394018, 0.08036454616199999, 0.00003897172772057774, 0.00003897172772057774),
    @("FAPs", "Vcan", "Tlr2", "ECs", 3, 1, 103.676216, 311.028648, 0.9862688099613843, 0.9862688099613843, 2, 0.6666666666666666, 0.493155, 1.479465, 0.2262533155038342, 0.2262533155038342, 51.12844430148, 460.1559987133199, 0.2231465882317842, 0.2231465882317842),
    @("FAPs", "Vcan", "Tlr2", "FAPs", 3, 1, 103.676216, 311.028648, 0.9862688099613843, 0.9862688099613843, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.7564224123165462, 0.7564224123165462, 170.9354008377493, 1538.418607539744, 0.7460358324235595, 0.7460358324235595),
    @("FAPs", "Vcan", "Tlr2", "MuSCs", 3, 1, 103.676216, 311.028648, 0.9862688099613843, 0.9862688099613843, 1, 0.3333333333333333, 0.037761, 0.113283, 0.01732427217961956, 0.01732427217961956, 3.914917592375999, 35.234258331384, 0.0170863893060405, 0.0170863893060405),
    @("MuSCs", "Vcan", "Tlr2", "ECs", 3, 1, 1.206946333333333, 3.620839, 0.01148164516212593, 0.01148164516212593, 2, 0.6666666666666666, 0.493155, 1.479465, 0.2262533155038342, 0.2262533155038342, 0.5952116190150001, 5.356904571135001, 0.00259776028536955, 0.002597760285369551),
    @("MuSCs", "Vcan", "Tlr2", "FAPs", 3, 1, 1.206946333333333, 3.620839, 0.01148164516212593, 0.01148164516212593, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.7564224123165462, 0.7564224123165462, 1.989943916143555, 17.909495245292, 0.008684973730897898, 0.0086849737308979),
    @("MuSCs", "Vcan", "Tlr2", "MuSCs", 3, 1, 1.206946333333333, 3.620839, 0.01148164516212593, 0.01148164516212593, 1, 0.3333333333333333, 0.037761, 0.113283, 0.01732427217961956, 0.01732427217961956, 0.045575500493, 0.410179504437, 0.0001989111458584818, 0.0001989111458584818)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "{0}{1}" -f $cols[$j], $r
        $ws.Range($addr).Value = $row[$j]
    }
}
